$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates: only the cells that actually change values are listed.
# Column D holds "Price" values that look numeric (e.g. "1.003", "93.20").
# They must stay TEXT (as in the source file) so formatting like trailing
# zeros and thousand-separator dots ("24.494.93") is preserved exactly -
# assigning such a string straight to .Value would make Excel coerce it
# into a real number and silently drop the formatting. Forcing the cell's
# NumberFormat to "@" (Text) first guarantees the literal string is kept.
$updates = @(
    @{ Row = 2;  D = "24.494.93";   E = "  -1.60%  " }
    @{ Row = 3;  D = "1.670.35";    E = "  -1.88%  " }
    @{ Row = 4;  D = "1.003";       E = "  -0.12%  " }
    @{ Row = 5;  D = "313.26";      E = "  -0.70%  " }
    @{ Row = 6;  E = "  -0.10%  " }
    @{ Row = 7;  D = "0.3905";      E = "  -3.44%  " }
    @{ Row = 8;  D = "0.3929";      E = "  -3.19%  " }
    @{ Row = 9;  D = "1.004";       E = "  +0.08%  " }
    @{ Row = 10; D = "51.78";       E = "  -3.55%  " }
    @{ Row = 11; D = "1.399";       E = "  -4.71%  " }
    @{ Row = 12; D = "0.08631";     E = "  -2.00%  " }
    @{ Row = 13; D = "25.27";       E = "  -2.07%  " }
    @{ Row = 14; D = "7.281";       E = "  -3.23%  " }
    @{ Row = 15; D = "0.00001319";  E = "  -2.38%  " }
    @{ Row = 16; D = "7.703";       E = "  -4.26%  " }
    @{ Row = 17; D = "1.673.36";    E = "  -3.35%  " }
    @{ Row = 18; D = "93.20";       E = "  -3.37%  " }
    @{ Row = 19; D = "0.07057";     E = "  -1.54%  " }
    @{ Row = 20; E = "  -1.28%  " }
    @{ Row = 21; D = "7.054";       E = "  -2.46%  " }
    @{ Row = 22; E = "  +0.00%  " }
    @{ Row = 23; D = "13.95";       E = "  -4.36%  " }
    @{ Row = 24; D = "24.475.88";   E = "  -1.69%  " }
    @{ Row = 25; D = "2.376";       E = "  +2.18%  " }
    @{ Row = 26; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "23.17";  E = "  +0.62%  " }
    @{ Row = 27; B = "LidoDAOToken";    C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo";       D = "2.735"; E = "  -5.28%  " }
    @{ Row = 28; D = "160.86";      E = "  -2.44%  " }
    @{ Row = 29; D = "5.846";       E = "  -14.47%  " }
    @{ Row = 30; D = "147.60";      E = "  +1.73%  " }
    @{ Row = 31; D = "8.268";       E = "  +0.42%  " }
    @{ Row = 32; D = "2.512";       E = "  +10.72%  " }
    @{ Row = 33; D = "1.857.54";    E = "  -3.09%  " }
    @{ Row = 34; D = "0.08344";     E = "  -5.24%  " }
    @{ Row = 35; D = "6.975";       E = "  -4.62%  " }
    @{ Row = 36; D = "0.03018";     E = "  -5.69%  " }
    @{ Row = 37; D = "0.2806";      E = "  -1.35%  " }
    @{ Row = 38; D = "0.9796";      E = "  -3.36%  " }
    @{ Row = 39; D = "0.09442";     E = "  +0.46%  " }
    @{ Row = 40; D = "1.526";       E = "  +3.93%  " }
    @{ Row = 41; D = "10.32";       E = "  -4.24%  " }
    @{ Row = 42; D = "0.7879";      E = "  -6.75%  " }
    @{ Row = 43; D = "13.55";       E = "  -3.72%  " }
    @{ Row = 44; E = "  -9.31%  " }
    @{ Row = 45; D = "0.7108";      E = "  -4.29%  " }
    @{ Row = 46; D = "2.546";       E = "  -6.14%  " }
    @{ Row = 47; D = "4.174";       E = "  -1.57%  " }
    @{ Row = 48; E = "  -0.09%  " }
    @{ Row = 49; D = "0.08576";     E = "  +2.75%  " }
    @{ Row = 50; D = "1.320";       E = "  -5.12%  " }
    @{ Row = 51; D = "137.49";      E = "  -3.24%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $dVal = $u.D
        # Only force Text format when the literal would otherwise be
        # auto-coerced into a number (losing e.g. trailing zeros / the
        # "thousands.decimal" look of values like "24.494.93"). Values
        # with two dots are never parsed as numbers by Excel, so they can
        # be written straight through and keep the cell's original (no
        # explicit style) formatting, just like in the source file.
        $looksNumeric = $dVal -match '^[+-]?[0-9]*\.?[0-9]+$'
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $dVal
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
